$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 5163.2
$ws.Range("I2").Value = 10141.6
$ws.Range("J2").Value = 184.8
$ws.Range("K2").Value = 10141.6
$ws.Range("L2").Value = 184.8
$ws.Range("M2").Value = -10028.6
$ws.Range("N2").Value = -410.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2665
$ws.Range("I100").Value = 1995
$ws.Range("K100").Value = 1995
$ws.Range("M100").Value = -1454

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 47651.715
$ws.Range("J111").Value = 59765.5
$ws.Range("L111").Value = 179296.5
$ws.Range("N111").Value = -185430.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3742.51
$ws.Range("I137").Value = 2237.25
$ws.Range("J137").Value = 3873.402
$ws.Range("K137").Value = 6711.75
$ws.Range("L137").Value = 11620.206
$ws.Range("M137").Value = -4161.75
$ws.Range("N137").Value = -16720.206

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5857.759
$ws.Range("I138").Value = 5991.8
$ws.Range("K138").Value = 17975.4
$ws.Range("M138").Value = -12835.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1854351.6
$ws.Range("I45").Value = 2779527.5
$ws.Range("K45").Value = 2779527.5
$ws.Range("M45").Value = -2779150.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2170.5588
$ws.Range("I61").Value = 1736.9259
$ws.Range("J61").Value = 3843.1428
$ws.Range("K61").Value = 1736.9259
$ws.Range("L61").Value = 3843.1428
$ws.Range("M61").Value = -1524.9259
$ws.Range("N61").Value = -4267.1428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 32072.428
$ws.Range("I102").Value = 35003.332
$ws.Range("K102").Value = 35003.332
$ws.Range("M102").Value = -33381.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 250039250
$ws.Range("I110").Value = 500002500
$ws.Range("J110").Value = 76000
$ws.Range("K110").Value = 500002500
$ws.Range("L110").Value = 76000
$ws.Range("M110").Value = -500000455
$ws.Range("N110").Value = -80090

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2170.5588
$ws.Range("I136").Value = 1736.9259
$ws.Range("J136").Value = 3843.1428
$ws.Range("K136").Value = 5210.7777
$ws.Range("L136").Value = 11529.4284
$ws.Range("M136").Value = -2660.7777
$ws.Range("N136").Value = -16629.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 33166.668
$ws.Range("J88").Value = 33166.668
$ws.Range("L88").Value = 33166.668
$ws.Range("N88").Value = -33978.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H91").Value = 33166.668
$ws.Range("J91").Value = 33166.668
$ws.Range("L91").Value = 33166.668
$ws.Range("N91").Value = -35974.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 27782166
$ws.Range("I94").Value = 50003200
$ws.Range("J94").Value = 5875.25
$ws.Range("K94").Value = 50003200
$ws.Range("L94").Value = 5875.25
$ws.Range("M94").Value = -50002749
$ws.Range("N94").Value = -6777.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 100073000
$ws.Range("I107").Value = 80000
$ws.Range("J107").Value = 166734990
$ws.Range("K107").Value = 80000
$ws.Range("L107").Value = 166734990
$ws.Range("M107").Value = -78080
$ws.Range("N107").Value = -166738830

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H118").Value = 125000
$ws.Range("J118").Value = 125000
$ws.Range("L118").Value = 125000
$ws.Range("N118").Value = -128314

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3139.8147
$ws.Range("I134").Value = 2791.125
$ws.Range("K134").Value = 8373.375
$ws.Range("M134").Value = -5838.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 2000
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -2574

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4115.7915
$ws.Range("I31").Value = 1998.4
$ws.Range("J31").Value = 4273.806
$ws.Range("K31").Value = 1998.4
$ws.Range("L31").Value = 4273.806
$ws.Range("M31").Value = -1703.4
$ws.Range("N31").Value = -4863.806

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4115.7915
$ws.Range("I34").Value = 1998.4
$ws.Range("J34").Value = 4273.806
$ws.Range("K34").Value = 1998.4
$ws.Range("L34").Value = 4273.806
$ws.Range("M34").Value = -1796.4
$ws.Range("N34").Value = -4677.806

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 30271
$ws.Range("J64").Value = 30271
$ws.Range("L64").Value = 30271
$ws.Range("N64").Value = -30767

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H67").Value = 30271
$ws.Range("J67").Value = 30271
$ws.Range("L67").Value = 30271
$ws.Range("N67").Value = -31987

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2110.56
$ws.Range("I107").Value = 1344.381
$ws.Range("K107").Value = 1344.381
$ws.Range("M107").Value = 575.6189999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1882.0883
$ws.Range("I134").Value = 1403.4073
$ws.Range("K134").Value = 4210.2219
$ws.Range("M134").Value = -1675.2219

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1590
$ws.Range("I97").Value = 1760.6
$ws.Range("J97").Value = 1468.1428
$ws.Range("K97").Value = 5281.799999999999
$ws.Range("L97").Value = 4404.428400000001
$ws.Range("M97").Value = -4785.799999999999
$ws.Range("N97").Value = -5396.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 7277
$ws.Range("I98").Value = 15053
$ws.Range("J98").Value = 2093
$ws.Range("K98").Value = 45159
$ws.Range("L98").Value = 6279
$ws.Range("M98").Value = -43661
$ws.Range("N98").Value = -9275

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1838.2222
$ws.Range("I107").Value = 1520.2858
$ws.Range("K107").Value = 4560.857400000001
$ws.Range("M107").Value = -2640.857400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1886.3334
$ws.Range("I114").Value = 829.5
$ws.Range("K114").Value = 2488.5
$ws.Range("M114").Value = 765.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 336.3913
$ws.Range("J122").Value = 363
$ws.Range("L122").Value = 3267
$ws.Range("N122").Value = -8167

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1373.0834
$ws.Range("I97").Value = 1393.8518
$ws.Range("J97").Value = 1310.7778
$ws.Range("K97").Value = 1393.8518
$ws.Range("L97").Value = 1310.7778
$ws.Range("M97").Value = -897.8517999999999
$ws.Range("N97").Value = -2302.7778

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2133.5
$ws.Range("I107").Value = 1505.6666
$ws.Range("J107").Value = 3075.25
$ws.Range("K107").Value = 1505.6666
$ws.Range("L107").Value = 3075.25
$ws.Range("M107").Value = 414.3334
$ws.Range("N107").Value = -6915.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3145.318
$ws.Range("I113").Value = 1919.7
$ws.Range("K113").Value = 1919.7
$ws.Range("M113").Value = 250.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 241502.95
$ws.Range("I132").Value = 325326.53
$ws.Range("J132").Value = 5272.909
$ws.Range("K132").Value = 975979.5900000001
$ws.Range("L132").Value = 15818.727
$ws.Range("M132").Value = -973449.5900000001
$ws.Range("N132").Value = -20878.727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 54999.832
$ws.Range("J141").Value = 93333
$ws.Range("L141").Value = 93333
$ws.Range("N141").Value = -103693

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 5000
$ws.Range("I93").Value = 5000
$ws.Range("K93").Value = 5000
$ws.Range("M93").Value = -3752

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4614.7075
$ws.Range("I132").Value = 2684.4583
$ws.Range("J132").Value = 7339.7646
$ws.Range("K132").Value = 8053.374899999999
$ws.Range("L132").Value = 22019.2938
$ws.Range("M132").Value = -5523.374899999999
$ws.Range("N132").Value = -27079.2938

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3805.0908
$ws.Range("I136").Value = 3595.9268
$ws.Range("K136").Value = 10787.7804
$ws.Range("M136").Value = -8237.7804

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 35000
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 35000
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 850.3333
$ws.Range("I107").Value = 874
$ws.Range("J107").Value = 803
$ws.Range("K107").Value = 2622
$ws.Range("L107").Value = 2409
$ws.Range("M107").Value = -702
$ws.Range("N107").Value = -6249

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2593.2
$ws.Range("I136").Value = 1614.8966
$ws.Range("K136").Value = 4844.6898
$ws.Range("M136").Value = -2294.6898
